$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "Hóa đơn" "Màn hình sản phẩm"
Replace-Text "Màn hình hóa đơn" "Màn hình hiển thị các sản phẩm"
Replace-Text "Thống kê" "Màn hình phân loại"
Replace-Text "Màn hình thống kê" "Màn hình phân loại các nhóm sản phẩm"
